$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.315.24"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "2.516.68"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "537.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.74%  "
$ws.Range("E7").Value = "  +0.32%  "
$ws.Range("D9").Value = "2.519.96"
$ws.Range("E9").Value = "  -0.82%  "
$ws.Range("E10").Value = "  +0.48%  "
$ws.Range("E11").Value = "  +1.46%  "
$ws.Range("E12").Value = "  -3.25%  "
$ws.Range("E13").Value = "  -0.60%  "
$ws.Range("D14").Value = "2.962.61"
$ws.Range("E14").Value = "  -0.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.48"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.67%  "
$ws.Range("D16").Value = "59.214.41"
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").Value = "2.515.27"
$ws.Range("E18").Value = "  -0.90%  "
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("E20").Value = "  +0.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "324.83"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.82"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.46"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.32%  "
$ws.Range("E25").Value = "  -0.88%  "
$ws.Range("E26").Value = "  +1.76%  "
$ws.Range("E27").Value = "  +0.55%  "
$ws.Range("E28").Value = "  -2.19%  "
$ws.Range("E29").Value = "  +3.86%  "
$ws.Range("D30").Value = "0.0₃0777"
$ws.Range("E30").Value = "  -0.57%  "
$ws.Range("E31").Value = "  -1.83%  "
$ws.Range("E32").Value = "  +5.07%  "
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("E34").Value = "  -2.51%  "
$ws.Range("E35").Value = "  -6.96%  "
$ws.Range("E36").Value = "  -1.13%  "
$ws.Range("E37").Value = "  -2.26%  "
$ws.Range("E38").Value = "  -1.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.88"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.22%  "
$ws.Range("E40").Value = "  -0.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.815"
$ws.Range("D41").Style = "Normal"
$ws.Range("E42").Value = "  -6.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "279.87"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.92%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.86"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.598"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0929"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.65%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "122.88"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.46%  "
$ws.Range("E49").Value = "  -0.26%  "
$ws.Range("E50").Value = "  -1.82%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.84"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.70%  "
